$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "29.464.51"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.849.47"
$ws.Range("E3").Value = "  -0.56%  "
Set-TextValue "D4" "0.9989"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "241.85"
$ws.Range("E5").Value = "  -1.22%  "
Set-TextValue "D6" "0.6264"
$ws.Range("E6").Value = "  -2.67%  "
Set-TextValue "D7" "0.9998"
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue "D8" "48.29"
$ws.Range("E8").Value = "  +1.17%  "
Set-TextValue "D9" "0.07524"
$ws.Range("E9").Value = "  -0.36%  "
Set-TextValue "D10" "0.2971"
$ws.Range("E10").Value = "  -0.08%  "
Set-TextValue "D11" "24.26"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.955.70"
$ws.Range("E12").Value = "  +5.13%  "
Set-TextValue "D13" "0.07701"
$ws.Range("E13").Value = "  +0.19%  "
Set-TextValue "D14" "4.994"
$ws.Range("E14").Value = "  -1.15%  "
Set-TextValue "D15" "0.6844"
$ws.Range("E15").Value = "  -1.25%  "
Set-TextValue "D16" "83.65"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "2.240.04"
$ws.Range("E17").Value = "  +6.08%  "
Set-TextValue "D18" "0.000009737"
$ws.Range("E18").Value = "  -1.62%  "
Set-TextValue "D19" "6.217"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "29.599.37"
$ws.Range("E20").Value = "  -0.67%  "
Set-TextValue "D21" "233.40"
$ws.Range("E21").Value = "  -1.65%  "
Set-TextValue "D22" "12.47"
$ws.Range("E22").Value = "  -1.50%  "
Set-TextValue "D23" "0.9997"
$ws.Range("E23").Value = "  -0.04%  "
Set-TextValue "D24" "7.562"
$ws.Range("E24").Value = "  +0.51%  "
Set-TextValue "D25" "0.9999"
$ws.Range("E25").Value = "  -0.09%  "
Set-TextValue "D26" "154.96"
$ws.Range("E26").Value = "  -2.42%  "
Set-TextValue "D27" "0.1390"
$ws.Range("E27").Value = "  -2.34%  "
Set-TextValue "D28" "8.427"
$ws.Range("E28").Value = "  -1.56%  "
Set-TextValue "D29" "17.69"
$ws.Range("E29").Value = "  -1.30%  "
Set-TextValue "D30" "1.480"
$ws.Range("E30").Value = "  -1.30%  "
Set-TextValue "D31" "0.05856"
$ws.Range("E31").Value = "  -6.27%  "
Set-TextValue "D32" "1.263"
$ws.Range("E32").Value = "  -2.20%  "
Set-TextValue "D33" "4.091"
$ws.Range("E33").Value = "  -1.52%  "
Set-TextValue "D34" "4.035"
$ws.Range("E34").Value = "  -1.72%  "
Set-TextValue "D35" "1.896"
Set-TextValue "D36" "1.168"
$ws.Range("E36").Value = "  -0.53%  "
Set-TextValue "D37" "0.7193"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "1.238.98"
Set-TextValue "D40" "2.792"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +0.03%  "
Set-TextValue "D42" "0.9065"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").Value = "2.137.04"
$ws.Range("E44").Value = "  +5.63%  "
Set-TextValue "D45" "0.9998"
$ws.Range("E45").Value = "  -0.06%  "
Set-TextValue "D46" "102.05"
$ws.Range("E46").Value = "  +0.08%  "
Set-TextValue "D47" "66.99"
$ws.Range("E47").Value = "  -0.37%  "
Set-TextValue "D48" "7.300"
$ws.Range("E48").Value = "  +8.39%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000119"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "1.724"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "9.158"
$ws.Range("E51").Value = "  -0.81%  "
